$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = @(
    "Mon Oct 09 22:39:47 EDT 2023",
    "Mon Oct 09 22:40:00 EDT 2023",
    "Mon Oct 09 22:40:13 EDT 2023",
    "Mon Oct 09 22:40:26 EDT 2023",
    "Mon Oct 09 22:40:39 EDT 2023",
    "Mon Oct 09 22:40:51 EDT 2023",
    "Mon Oct 09 22:41:04 EDT 2023",
    "Mon Oct 09 22:41:17 EDT 2023",
    "Mon Oct 09 22:41:30 EDT 2023",
    "Mon Oct 09 22:41:43 EDT 2023",
    "Mon Oct 09 22:41:56 EDT 2023",
    "Mon Oct 09 22:42:08 EDT 2023",
    "Mon Oct 09 22:42:21 EDT 2023",
    "Mon Oct 09 22:42:34 EDT 2023",
    "Mon Oct 09 22:42:47 EDT 2023",
    "Mon Oct 09 22:43:00 EDT 2023",
    "Mon Oct 09 22:43:12 EDT 2023",
    "Mon Oct 09 22:43:25 EDT 2023",
    "Mon Oct 09 22:43:38 EDT 2023",
    "Mon Oct 09 22:43:51 EDT 2023",
    "Mon Oct 09 22:44:03 EDT 2023",
    "Mon Oct 09 22:44:16 EDT 2023",
    "Mon Oct 09 22:44:29 EDT 2023",
    "Mon Oct 09 22:44:42 EDT 2023",
    "Mon Oct 09 22:44:55 EDT 2023",
    "Mon Oct 09 22:45:08 EDT 2023",
    "Mon Oct 09 22:45:20 EDT 2023",
    "Mon Oct 09 22:45:33 EDT 2023",
    "Mon Oct 09 22:45:47 EDT 2023",
    "Mon Oct 09 22:45:59 EDT 2023",
    "Mon Oct 09 22:46:12 EDT 2023",
    "Mon Oct 09 22:46:25 EDT 2023"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $timestamps[$i]
}
